$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format for numeric-looking price cells so exact formatting (trailing zeros, etc.) is preserved
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range('D2').Value = '60.350.33'
$ws.Range('E2').Value = '  +1.88%  '
$ws.Range('D3').Value = '2.596.14'
$ws.Range('E3').Value = '  +0.55%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '569.35'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').Value = '141.90'
$ws.Range('E6').Value = '  -0.52%  '
$ws.Range('D7').Value = '0.997'
$ws.Range('E7').Value = '  -0.20%  '
$ws.Range('D8').Value = '0.599'
$ws.Range('E8').Value = '  +0.45%  '
$ws.Range('D9').Value = '2.617.67'
$ws.Range('E9').Value = '  +1.07%  '
$ws.Range('D10').Value = '6.57'
$ws.Range('E10').Value = '  -0.83%  '
$ws.Range('E11').Value = '  +1.92%  '
$ws.Range('D12').Value = '0.368'
$ws.Range('E12').Value = '  +4.76%  '
$ws.Range('E13').Value = '  -6.28%  '
$ws.Range('D14').Value = '3.056.38'
$ws.Range('E14').Value = '  +0.63%  '
$ws.Range('D15').Value = '60.342.32'
$ws.Range('E15').Value = '  +1.94%  '
$ws.Range('D16').Value = '23.38'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '0.0000141'
$ws.Range('E17').Value = '  +2.84%  '
$ws.Range('D18').Value = '2.607.50'
$ws.Range('E18').Value = '  +0.87%  '
$ws.Range('D19').Value = '11.30'
$ws.Range('E19').Value = '  +9.04%  '
$ws.Range('D20').Value = '4.66'
$ws.Range('E20').Value = '  +1.90%  '
$ws.Range('D21').Value = '346.65'
$ws.Range('E21').Value = '  +2.77%  '
$ws.Range('D22').Value = '6.97'
$ws.Range('E22').Value = '  +9.02%  '
$ws.Range('D23').Value = '0.998'
$ws.Range('E23').Value = '  -0.28%  '
$ws.Range('D24').Value = '0.535'
$ws.Range('E24').Value = '  +14.66%  '
$ws.Range('D25').Value = '63.16'
$ws.Range('E25').Value = '  -1.56%  '
$ws.Range('E26').Value = '  -0.30%  '
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('E28').Value = '  +4.93%  '
$ws.Range('D29').Value = '0.0₃0785'
$ws.Range('E29').Value = '  +1.56%  '
$ws.Range('E30').Value = '  +9.60%  '
$ws.Range('D31').Value = '0.997'
$ws.Range('E31').Value = '  -0.08%  '
$ws.Range('D32').Value = '6.32'
$ws.Range('E32').Value = '  +3.61%  '
$ws.Range('D33').Value = '161.14'
$ws.Range('E33').Value = '  +0.80%  '
$ws.Range('D34').Value = '19.45'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('E35').Value = '  +4.63%  '
$ws.Range('D36').Value = '0.962'
$ws.Range('E36').Value = '  +9.62%  '
$ws.Range('E37').Value = '  +4.25%  '
$ws.Range('E38').Value = '  +8.40%  '
$ws.Range('D39').Value = '37.81'
$ws.Range('E39').Value = '  +1.06%  '
$ws.Range('E40').Value = '  +4.03%  '
$ws.Range('D41').Value = '0.848'
$ws.Range('E41').Value = '  -2.77%  '
$ws.Range('D42').Value = '294.24'
$ws.Range('E42').Value = '  -0.12%  '
$ws.Range('D43').Value = '137.87'
$ws.Range('E43').Value = '  +4.61%  '
$ws.Range('D44').Value = '0.996'
$ws.Range('E44').Value = '  -0.30%  '
$ws.Range('D45').Value = '0.608'
$ws.Range('E45').Value = '  +2.04%  '
$ws.Range('D46').Value = '0.0981'
$ws.Range('E46').Value = '  +0.91%  '
$ws.Range('D47').Value = '19.64'
$ws.Range('E47').Value = '  +3.17%  '
$ws.Range('D48').Value = '0.0546'
$ws.Range('E48').Value = '  +1.94%  '
$ws.Range('D49').Value = '19.93'
$ws.Range('E49').Value = '  +7.14%  '
$ws.Range('D50').Value = '0.0240'
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('B51').Value = 'RenderToken'
$ws.Range('C51').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D51').Value = '4.89'
$ws.Range('E51').Value = '  +8.72%  '
